$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 381.76923
$ws.Range("I33").Value = 259.8889
$ws.Range("J33").Value = 656
$ws.Range("K33").Value = 259.8889
$ws.Range("L33").Value = 656
$ws.Range("M33").Value = -30.88889999999998
$ws.Range("N33").Value = -1114

$ws.Range("H62").Value = 5532.731
$ws.Range("I62").Value = 5116.952
$ws.Range("J62").Value = 7279
$ws.Range("K62").Value = 5116.952
$ws.Range("L62").Value = 7279
$ws.Range("M62").Value = -4492.952
$ws.Range("N62").Value = -8527

$ws.Range("H65").Value = 5532.731
$ws.Range("I65").Value = 5116.952
$ws.Range("J65").Value = 7279
$ws.Range("K65").Value = 25584.76
$ws.Range("L65").Value = 36395
$ws.Range("M65").Value = -22464.76
$ws.Range("N65").Value = -42635

$ws.Range("H103").Value = 843.2308
$ws.Range("I103").Value = 643.8
$ws.Range("J103").Value = 967.875
$ws.Range("K103").Value = 1931.4
$ws.Range("L103").Value = 2903.625
$ws.Range("M103").Value = -1345.4
$ws.Range("N103").Value = -4075.625

$ws.Range("H107").Value = 1332.3125
$ws.Range("J107").Value = 1587.5
$ws.Range("L107").Value = 1587.5
$ws.Range("N107").Value = -5427.5

$ws.Range("H113").Value = 4498.5
$ws.Range("I113").Value = 2998.5
$ws.Range("J113").Value = 5248.5
$ws.Range("K113").Value = 2998.5
$ws.Range("L113").Value = 5248.5
$ws.Range("M113").Value = 255.5
$ws.Range("N113").Value = -11756.5

$ws.Range("H134").Value = 50000
$ws.Range("J134").Value = 50000
$ws.Range("L134").Value = 50000
$ws.Range("N134").Value = -60140

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 3533.5
$ws.Range("I45").Value = 1867.25
$ws.Range("K45").Value = 1867.25
$ws.Range("M45").Value = -1490.25

$ws.Range("H61").Value = 4938.48
$ws.Range("I61").Value = 612.2857
$ws.Range("J61").Value = 27651
$ws.Range("K61").Value = 612.2857
$ws.Range("L61").Value = 27651
$ws.Range("M61").Value = -400.2857
$ws.Range("N61").Value = -28075

$ws.Range("H110").Value = 545.4286
$ws.Range("I110").Value = 553.6667
$ws.Range("K110").Value = 553.6667
$ws.Range("M110").Value = 1491.3333

$ws.Range("H132").Value = 1391.814
$ws.Range("I132").Value = 1072.9697
$ws.Range("K132").Value = 3218.9091
$ws.Range("M132").Value = -688.9091000000003

$ws.Range("H136").Value = 4938.48
$ws.Range("I136").Value = 612.2857
$ws.Range("J136").Value = 27651
$ws.Range("K136").Value = 1836.8571
$ws.Range("L136").Value = 82953
$ws.Range("M136").Value = 713.1428999999998
$ws.Range("N136").Value = -88053

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 32714.2
$ws.Range("J20").Value = 2532.8
$ws.Range("L20").Value = 2532.8
$ws.Range("N20").Value = -3026.8

$ws.Range("H94").Value = 960.67566
$ws.Range("I94").Value = 835.25
$ws.Range("J94").Value = 1192.2307
$ws.Range("K94").Value = 835.25
$ws.Range("L94").Value = 1192.2307
$ws.Range("M94").Value = -384.25
$ws.Range("N94").Value = -2094.2307

$ws.Range("H105").Value = 2955.5557
$ws.Range("I105").Value = 1992.2727
$ws.Range("J105").Value = 4469.2856
$ws.Range("K105").Value = 1992.2727
$ws.Range("L105").Value = 4469.2856
$ws.Range("M105").Value = -245.2727
$ws.Range("N105").Value = -7963.2856

$ws.Range("H134").Value = 14690.619
$ws.Range("I134").Value = 17469.312
$ws.Range("K134").Value = 52407.936
$ws.Range("M134").Value = -49872.936

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 627
$ws.Range("I22").Value = 412.2857
$ws.Range("J22").Value = 841.7143
$ws.Range("K22").Value = 412.2857
$ws.Range("L22").Value = 841.7143
$ws.Range("M22").Value = -62.28570000000002
$ws.Range("N22").Value = -1541.7143

$ws.Range("H31").Value = 25005754
$ws.Range("I31").Value = 100000000
$ws.Range("K31").Value = 100000000
$ws.Range("M31").Value = -99999705

$ws.Range("H34").Value = 25005754
$ws.Range("I34").Value = 100000000
$ws.Range("K34").Value = 100000000
$ws.Range("M34").Value = -99999798

$ws.Range("H58").Value = 1623.9697
$ws.Range("I58").Value = 1237.2084
$ws.Range("K58").Value = 1237.2084
$ws.Range("M58").Value = -1034.2084

$ws.Range("H99").Value = 9207.25
$ws.Range("I99").Value = 5685.5
$ws.Range("K99").Value = 5685.5
$ws.Range("M99").Value = -4187.5

$ws.Range("H107").Value = 616.3929000000001
$ws.Range("I107").Value = 509.83334
$ws.Range("J107").Value = 808.2
$ws.Range("K107").Value = 509.83334
$ws.Range("L107").Value = 808.2
$ws.Range("M107").Value = 1410.16666
$ws.Range("N107").Value = -4648.2

$ws.Range("H126").Value = 9207.25
$ws.Range("I126").Value = 5685.5
$ws.Range("K126").Value = 17056.5
$ws.Range("M126").Value = -14586.5

$ws.Range("H136").Value = 1623.9697
$ws.Range("I136").Value = 1237.2084
$ws.Range("K136").Value = 3711.6252
$ws.Range("M136").Value = -1161.6252

$ws.Range("H141").Value = 66690.836
$ws.Range("J141").Value = 66690.836
$ws.Range("L141").Value = 66690.836
$ws.Range("N141").Value = -77050.836

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1089.4
$ws.Range("I5").Value = 861.8333
$ws.Range("J5").Value = 1999.6666
$ws.Range("K5").Value = 2585.4999
$ws.Range("L5").Value = 5998.9998
$ws.Range("M5").Value = -2473.4999
$ws.Range("N5").Value = -6222.9998

$ws.Range("H37").Value = 41972.555
$ws.Range("J37").Value = 41972.555
$ws.Range("L37").Value = 125917.665
$ws.Range("N37").Value = -126141.665

$ws.Range("H94").Value = 2875
$ws.Range("I94").Value = 750
$ws.Range("K94").Value = 2250
$ws.Range("M94").Value = -1574

$ws.Range("H122").Value = 887.2105
$ws.Range("I122").Value = 586.3333
$ws.Range("J122").Value = 943.625
$ws.Range("K122").Value = 5276.9997
$ws.Range("L122").Value = 8492.625
$ws.Range("M122").Value = -2826.9997
$ws.Range("N122").Value = -13392.625

$ws.Range("H135").Value = 1089.4
$ws.Range("I135").Value = 861.8333
$ws.Range("J135").Value = 1999.6666
$ws.Range("K135").Value = 7756.4997
$ws.Range("L135").Value = 17996.9994
$ws.Range("M135").Value = -5221.4997
$ws.Range("N135").Value = -23066.9994

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H58").Value = 35199.9
$ws.Range("I58").Value = 15999.5
$ws.Range("K58").Value = 15999.5
$ws.Range("M58").Value = -15722.5

$ws.Range("H122").Value = 3898.9443
$ws.Range("I122").Value = 3754.5833
$ws.Range("K122").Value = 11263.7499
$ws.Range("M122").Value = -8813.749899999999

$ws.Range("H132").Value = 0
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 0
$ws.Range("M132").ClearContents()
$ws.Range("N132").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3541.1667
$ws.Range("I7").Value = 3861.75
$ws.Range("J7").Value = 2900
$ws.Range("K7").Value = 3861.75
$ws.Range("L7").Value = 2900
$ws.Range("M7").Value = -3749.75
$ws.Range("N7").Value = -3124

$ws.Range("H46").Value = 3786.7778
$ws.Range("I46").Value = 2317.3333
$ws.Range("K46").Value = 2317.3333
$ws.Range("M46").Value = -2129.3333

$ws.Range("H100").Value = 2827.9473
$ws.Range("I100").Value = 2732.8462
$ws.Range("J100").Value = 3034
$ws.Range("K100").Value = 2732.8462
$ws.Range("L100").Value = 3034
$ws.Range("M100").Value = -2191.8462
$ws.Range("N100").Value = -4116

$ws.Range("H126").Value = 3541.1667
$ws.Range("I126").Value = 3861.75
$ws.Range("J126").Value = 2900
$ws.Range("K126").Value = 11585.25
$ws.Range("L126").Value = 8700
$ws.Range("M126").Value = -9115.25
$ws.Range("N126").Value = -13640

$ws.Range("H136").Value = 4280.48
$ws.Range("I136").Value = 4102.25
$ws.Range("K136").Value = 12306.75
$ws.Range("M136").Value = -9756.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H59").Value = 30000
$ws.Range("J59").Value = 30000
$ws.Range("L59").Value = 30000
$ws.Range("N59").Value = -31476

$ws.Range("H126").Value = 4000
$ws.Range("I126").Value = 2000
$ws.Range("J126").Value = 6000
$ws.Range("K126").Value = 6000
$ws.Range("L126").Value = 18000
$ws.Range("M126").Value = -3530
$ws.Range("N126").Value = -22940

$ws.Range("H132").Value = 1428.0769
$ws.Range("I132").Value = 1041
$ws.Range("K132").Value = 3123
$ws.Range("M132").Value = -593
